# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
#
# This script appends the Week 17 per-play logs to the running season logs
# (YDS and ST sheets) and updates the season-to-date totals (OFF, DEF, ST,
# TURNS and PEN sheets) to include Week 17.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append week 17 per-play yardage logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$add_B2 = @"
 2 0 -1 4 10 5 3 4 7 2 7 7 9 2 2 2 10 4 4 1 4 -1 1 3 3 0
"@
$ws.Range("B2").Value2 = $ws.Range("B2").Value2 + $add_B2

$add_C2 = @"
 8 4 6 2 6 5 4 -1 8 5 1 2 10 8 2 3 3 5 2 -2 4 -3 -1 9 6 4 1 8 3 3 1 1 13 5 7 5 1 3 9 1
"@
$ws.Range("C2").Value2 = $ws.Range("C2").Value2 + $add_C2

$add_B3 = @"
 3 11 8 4 6 23 18 17 4 6 16 8 8 11 10 3 13
"@
$ws.Range("B3").Value2 = $ws.Range("B3").Value2 + $add_B3

$add_C3 = @"
 4 12 12 -4
"@
$ws.Range("C3").Value2 = $ws.Range("C3").Value2 + $add_C3

# ---------------------------------------------------------------------
# OFF sheet - season totals through week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value2 = 385
$ws.Range("D2").Value2 = 29
$ws.Range("F2").Value2 = 142
$ws.Range("G2").Value2 = 100
$ws.Range("I2").Value2 = 10
$ws.Range("J2").Value2 = 62
$ws.Range("L2").Value2 = 535
$ws.Range("M2").Value2 = 348
$ws.Range("O2").Value2 = 58
$ws.Range("P2").Value2 = 31
$ws.Range("Q2").Value2 = 1010

$ws.Range("B3").Value2 = 15
$ws.Range("C3").Value2 = 368
$ws.Range("E3").Value2 = 75
$ws.Range("F3").Value2 = 186
$ws.Range("G3").Value2 = 82
$ws.Range("H3").Value2 = 64
$ws.Range("I3").Value2 = 121
$ws.Range("J3").Value2 = 104
$ws.Range("N3").Value2 = 44

# ---------------------------------------------------------------------
# DEF sheet - season totals through week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value2 = 420
$ws.Range("E2").Value2 = 23
$ws.Range("F2").Value2 = 143
$ws.Range("G2").Value2 = 126
$ws.Range("I2").Value2 = 14
$ws.Range("J2").Value2 = 59
$ws.Range("L2").Value2 = 484
$ws.Range("M2").Value2 = 319
$ws.Range("O2").Value2 = 48
$ws.Range("Q2").Value2 = 973

$ws.Range("C3").Value2 = 323
$ws.Range("G3").Value2 = 66
$ws.Range("I3").Value2 = 112
$ws.Range("J3").Value2 = 128
$ws.Range("N3").Value2 = 46

# ---------------------------------------------------------------------
# ST sheet - season totals through week 17, plus appended week 17 logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value2 = 149
$ws.Range("D2").Value2 = 132
$ws.Range("F2").Value2 = 232
$ws.Range("G2").Value2 = 224
$ws.Range("J2").Value2 = 114
$ws.Range("K2").Value2 = 107
$ws.Range("L2").Value2 = 57
$ws.Range("M2").Value2 = 44

$ws.Range("B3").Value2 = 76

$add_B4 = @"
 66 62 63
"@
$ws.Range("B4").Value2 = $ws.Range("B4").Value2 + $add_B4

$add_B5 = @"
 26 24 3
"@
$ws.Range("B5").Value2 = $ws.Range("B5").Value2 + $add_B5

$add_B6 = @"
 16 6
"@
$ws.Range("B6").Value2 = $ws.Range("B6").Value2 + $add_B6

$add_D3 = @"
 42 42 56 52 35
"@
$ws.Range("D3").Value2 = $ws.Range("D3").Value2 + $add_D3

$add_D4 = @"
 15 0 17 0 0
"@
$ws.Range("D4").Value2 = $ws.Range("D4").Value2 + $add_D4

$add_D5 = @"
 1 0 0 0
"@
$ws.Range("D5").Value2 = $ws.Range("D5").Value2 + $add_D5

# ---------------------------------------------------------------------
# TURNS sheet - season totals through week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B2").Value2 = 15
$ws.Range("C2").Value2 = 11
$ws.Range("E2").Value2 = 19
$ws.Range("D3").Value2 = 14

# ---------------------------------------------------------------------
# PEN sheet - season totals through week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value2 = 40
$ws.Range("B3").Value2 = 24
$ws.Range("D3").Value2 = 7
